# Restructure the certificate template grid:
#  - Rename/re-map the header row (B..H) to the new column set.
#  - Rewrite row 2 and row 3 data to match the new columns.
#  - Drop the old trailing columns (I:L) that are no longer part of the grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
# A1 "Nombre Completo" is unchanged.
$ws.Range("B1").Value = "ID del Curso"
$ws.Range("C1").Value = "Identificación"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Teléfono"
$ws.Range("F1").Value = "Ubicación Física"
$ws.Range("G1").Value = "Estado"
$ws.Range("H1").Value = "Origen"

# ---- Row 2 (Juan Pérez) ----
# A2 "Juan Pérez" is unchanged.
$ws.Range("B2").Value = "NUT-1-2025"
$ws.Range("C2").Value = "'0801-1990-12345"
$ws.Range("D2").Value = "juan@ejemplo.com"
$ws.Range("E2").Value = "'50499887766"
$ws.Range("F2").Value = "Tomo 1 Caja 5"
$ws.Range("G2").Value = "en_archivo"
$ws.Range("H2").Value = "nuevo"

# ---- Row 3 (María García) ----
# A3 "María García" is unchanged.
$ws.Range("B3").Value = "SALUD-2024"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "maria@ejemplo.com"
$ws.Range("E3").Value = "'50488776655"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "entregado"
$ws.Range("H3").Value = "historico"

# ---- Drop old columns I:L (Email/Teléfono/Ubicación/Estado/Origen have
# ---- been absorbed into D:H above, so the tail columns are now empty) ----
$ws.Range("I1:L3").ClearContents()
